$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NATMI LR-pair data (Ngf-Ngfr), recomputed with updated TPM values.
# Columns A-T correspond to the header row 1 fields; 12 data rows (3 sending
# clusters x 4 target clusters) replace the previous 9 rows (3 x 3).
$data = @(
    @("ECs", "Ngf", "Ngfr", "ECs", 3, 1, 1.814087666666667, 5.442263, 0.1211063206477811, 0.1211063206477811, 2, 0.6666666666666666, 0.3922183333333333, 1.176655, 0.06257714082953221, 0.06257714082953222, 0.7115184411405555, 6.403665970264999, 0.007578487282522684, 0.007578487282522686),
    @("ECs", "Ngf", "Ngfr", "MuSCs", 3, 1, 1.814087666666667, 5.442263, 0.1211063206477811, 0.1211063206477811, 3, 1, 5.787580333333334, 17.362741, 0.9233893441524432, 0.9233893441524432, 10.49917810254255, 94.49260292288299, 0.1118282859956701, 0.1118282859956701),
    @("ECs", "Ngf", "Ngfr", "Neutrophils", 3, 1, 1.814087666666667, 5.442263, 0.1211063206477811, 0.1211063206477811, 1, 0.3333333333333333, 0.034619, 0.103857, 0.005523347213187152, 0.005523347213187152, 0.06280190093233333, 0.565217108391, 0.0006689122586492715, 0.0006689122586492715),
    @("ECs", "Ngf", "Ngfr", "Resolving-Mac", 3, 1, 1.814087666666667, 5.442263, 0.1211063206477811, 0.1211063206477811, 1, 0.3333333333333333, 0.05333966666666667, 0.160019, 0.008510167804837369, 0.008510167804837371, 0.0967628314441111, 0.870865482997, 0.001030635110939058, 0.001030635110939058),
    @("FAPs", "Ngf", "Ngfr", "ECs", 3, 1, 5.748364666666667, 17.245094, 0.3837539427192561, 0.3837539427192561, 2, 0.6666666666666666, 0.3922183333333333, 1.176655, 0.06257714082953221, 0.06257714082953222, 2.254614008952223, 20.29152608057, 0.02401422451743112, 0.02401422451743113),
    @("FAPs", "Ngf", "Ngfr", "MuSCs", 3, 1, 5.748364666666667, 17.245094, 0.3837539427192561, 0.3837539427192561, 3, 1, 5.787580333333334, 17.362741, 0.9233893441524432, 0.9233893441524432, 33.26912229362823, 299.422100642654, 0.3543543014834481, 0.3543543014834481),
    @("FAPs", "Ngf", "Ngfr", "Neutrophils", 3, 1, 5.748364666666667, 17.245094, 0.3837539427192561, 0.3837539427192561, 1, 0.3333333333333333, 0.034619, 0.103857, 0.005523347213187152, 0.005523347213187152, 0.1990026363953334, 1.791023727558, 0.002119606270067985, 0.002119606270067985),
    @("FAPs", "Ngf", "Ngfr", "Resolving-Mac", 3, 1, 5.748364666666667, 17.245094, 0.3837539427192561, 0.3837539427192561, 1, 0.3333333333333333, 0.05333966666666667, 0.160019, 0.008510167804837369, 0.008510167804837371, 0.3066158551984445, 2.759542696786, 0.003265810448308817, 0.003265810448308818),
    @("MuSCs", "Ngf", "Ngfr", "ECs", 3, 1, 7.416845666666667, 22.250537, 0.4951397366329628, 0.4951397366329628, 2, 0.6666666666666666, 0.3922183333333333, 1.176655, 0.06257714082953221, 0.06257714082953222, 2.909022845970556, 26.181205613735, 0.0309844290295784, 0.03098442902957841),
    @("MuSCs", "Ngf", "Ngfr", "MuSCs", 3, 1, 7.416845666666667, 22.250537, 0.4951397366329628, 0.4951397366329628, 3, 1, 5.787580333333334, 17.362741, 0.9233893441524432, 0.9233893441524432, 42.92559011576856, 386.330311041917, 0.4572067566733251, 0.4572067566733251),
    @("MuSCs", "Ngf", "Ngfr", "Neutrophils", 3, 1, 7.416845666666667, 22.250537, 0.4951397366329628, 0.4951397366329628, 1, 0.3333333333333333, 0.034619, 0.103857, 0.005523347213187152, 0.005523347213187152, 0.2567637801343334, 2.310874021209, 0.002734828684469896, 0.002734828684469896),
    @("MuSCs", "Ngf", "Ngfr", "Resolving-Mac", 3, 1, 7.416845666666667, 22.250537, 0.4951397366329628, 0.4951397366329628, 1, 0.3333333333333333, 0.05333966666666667, 0.160019, 0.008510167804837369, 0.008510167804837371, 0.3956120755781111, 3.560508680203, 0.004213722245589495, 0.004213722245589495)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
